$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 44456.5
$ws.Range("J3").Value = 44456.5
$ws.Range("L3").Value = 44456.5
$ws.Range("N3").Value = -44684.5
$ws.Range("H53").Value = 522.2222
$ws.Range("I53").Value = 139.66667
$ws.Range("K53").Value = 139.66667
$ws.Range("M53").Value = 497.33333
$ws.Range("H64").Value = 13892994
$ws.Range("I64").Value = 19234530
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 19234530
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -19234282
$ws.Range("N64").Value = -5496
$ws.Range("H67").Value = 13892994
$ws.Range("I67").Value = 19234530
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 19234530
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -19233672
$ws.Range("N67").Value = -6716
$ws.Range("H100").Value = 2015.875
$ws.Range("I100").Value = 1604.1538
$ws.Range("J100").Value = 3800
$ws.Range("K100").Value = 1604.1538
$ws.Range("L100").Value = 3800
$ws.Range("M100").Value = -1063.1538
$ws.Range("N100").Value = -4882
$ws.Range("H102").Value = 44456.5
$ws.Range("J102").Value = 44456.5
$ws.Range("L102").Value = 44456.5
$ws.Range("N102").Value = -50946.5
$ws.Range("H107").Value = 828.0769
$ws.Range("I107").Value = 907
$ws.Range("J107").Value = 650.5
$ws.Range("K107").Value = 907
$ws.Range("L107").Value = 650.5
$ws.Range("M107").Value = 1013
$ws.Range("N107").Value = -4490.5
$ws.Range("H113").Value = 2415.8
$ws.Range("I113").Value = 2639
$ws.Range("J113").Value = 2267
$ws.Range("K113").Value = 2639
$ws.Range("L113").Value = 2267
$ws.Range("M113").Value = 615
$ws.Range("N113").Value = -8775
$ws.Range("H116").Value = 4519.5386
$ws.Range("I116").Value = 4239.5654
$ws.Range("K116").Value = 4239.5654
$ws.Range("M116").Value = -797.5654000000004
$ws.Range("H125").Value = 11839127
$ws.Range("H132").Value = 5101.2173
$ws.Range("I132").Value = 5087.095
$ws.Range("J132").Value = 5249.5
$ws.Range("K132").Value = 15261.285
$ws.Range("L132").Value = 15748.5
$ws.Range("M132").Value = -12731.285
$ws.Range("N132").Value = -20808.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1899.5385
$ws.Range("I2").Value = 1862.125
$ws.Range("K2").Value = 1862.125
$ws.Range("M2").Value = -1749.125
$ws.Range("H32").Value = 5316.6387
$ws.Range("I32").Value = 2464.516
$ws.Range("K32").Value = 2464.516
$ws.Range("M32").Value = -2177.516
$ws.Range("H102").Value = 1142.6666
$ws.Range("I102").Value = 1034.1
$ws.Range("K102").Value = 1034.1
$ws.Range("M102").Value = 587.9000000000001
$ws.Range("H110").Value = 63160.25
$ws.Range("I110").Value = 63160.25
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 63160.25
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -61115.25
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 1899.5385
$ws.Range("I116").Value = 1862.125
$ws.Range("K116").Value = 1862.125
$ws.Range("M116").Value = 431.875
$ws.Range("H117").Value = 11999
$ws.Range("J117").Value = 11999
$ws.Range("L117").Value = 11999
$ws.Range("N117").Value = -21177

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1899.5385
$ws.Range("I3").Value = 1862.125
$ws.Range("K3").Value = 1862.125
$ws.Range("M3").Value = -1748.125
$ws.Range("H64").Value = 689.7273
$ws.Range("I64").Value = 477
$ws.Range("J64").Value = 737
$ws.Range("K64").Value = 477
$ws.Range("L64").Value = 737
$ws.Range("M64").Value = -252
$ws.Range("N64").Value = -1187
$ws.Range("H67").Value = 689.7273
$ws.Range("I67").Value = 477
$ws.Range("J67").Value = 737
$ws.Range("K67").Value = 477
$ws.Range("L67").Value = 737
$ws.Range("M67").Value = 303
$ws.Range("N67").Value = -2297
$ws.Range("H107").Value = 46164.695
$ws.Range("I107").Value = 2808.5908
$ws.Range("K107").Value = 2808.5908
$ws.Range("M107").Value = -888.5907999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 552
$ws.Range("I2").Value = 104
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 104
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = -1226
$ws.Range("H31").Value = 6098.5454
$ws.Range("J31").Value = 2956.2104
$ws.Range("L31").Value = 2956.2104
$ws.Range("N31").Value = -3546.2104
$ws.Range("H34").Value = 6098.5454
$ws.Range("J34").Value = 2956.2104
$ws.Range("L34").Value = 2956.2104
$ws.Range("N34").Value = -3360.2104
$ws.Range("H69").Value = 12599.667
$ws.Range("I69").Value = 12399.5
$ws.Range("K69").Value = 12399.5
$ws.Range("M69").Value = -11650.5
$ws.Range("H72").Value = 12599.667
$ws.Range("I72").Value = 12399.5
$ws.Range("K72").Value = 37198.5
$ws.Range("M72").Value = -33454.5
$ws.Range("H86").Value = 11404.462
$ws.Range("I86").Value = 9693
$ws.Range("J86").Value = 13401.167
$ws.Range("K86").Value = 9693
$ws.Range("L86").Value = 13401.167
$ws.Range("M86").Value = -8570
$ws.Range("N86").Value = -15647.167
$ws.Range("H89").Value = 11404.462
$ws.Range("I89").Value = 9693
$ws.Range("J89").Value = 13401.167
$ws.Range("K89").Value = 48465
$ws.Range("L89").Value = 67005.83499999999
$ws.Range("M89").Value = -42849
$ws.Range("N89").Value = -78237.83499999999
$ws.Range("H134").Value = 41672230
$ws.Range("I134").Value = 62504600
$ws.Range("J134").Value = 7497
$ws.Range("K134").Value = 187513800
$ws.Range("L134").Value = 22491
$ws.Range("M134").Value = -187511265
$ws.Range("N134").Value = -27561

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3587.1428
$ws.Range("J68").Value = 3865.4348
$ws.Range("L68").Value = 11596.3044
$ws.Range("N68").Value = -13218.3044
$ws.Range("H71").Value = 3587.1428
$ws.Range("J71").Value = 3865.4348
$ws.Range("L71").Value = 34788.9132
$ws.Range("N71").Value = -42900.9132
$ws.Range("H107").Value = 1936.238
$ws.Range("I107").Value = 676.7143
$ws.Range("J107").Value = 2566
$ws.Range("K107").Value = 2030.1429
$ws.Range("L107").Value = 7698
$ws.Range("M107").Value = -110.1428999999998
$ws.Range("N107").Value = -11538
$ws.Range("H122").Value = 813.4
$ws.Range("I122").Value = 760.8
$ws.Range("J122").Value = 918.6
$ws.Range("K122").Value = 6847.2
$ws.Range("L122").Value = 8267.4
$ws.Range("M122").Value = -4397.2
$ws.Range("N122").Value = -13167.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 29266
$ws.Range("J20").Value = 29899
$ws.Range("L20").Value = 29899
$ws.Range("N20").Value = -30389
$ws.Range("H102").Value = 4084.1667
$ws.Range("I102").Value = 4151
$ws.Range("K102").Value = 4151
$ws.Range("M102").Value = -2529
$ws.Range("H132").Value = 9617855
$ws.Range("I132").Value = 11366212
$ws.Range("J132").Value = 1891.5
$ws.Range("K132").Value = 34098636
$ws.Range("L132").Value = 5674.5
$ws.Range("M132").Value = -34096106
$ws.Range("N132").Value = -10734.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2876.5557
$ws.Range("I61").Value = 3155.1072
$ws.Range("J61").Value = 1901.625
$ws.Range("K61").Value = 3155.1072
$ws.Range("L61").Value = 1901.625
$ws.Range("M61").Value = -2953.1072
$ws.Range("N61").Value = -2305.625
$ws.Range("H68").Value = 3293973.5
$ws.Range("I68").Value = 4388214.5
$ws.Range("K68").Value = 4388214.5
$ws.Range("M68").Value = -4387465.5
$ws.Range("H71").Value = 3293973.5
$ws.Range("I71").Value = 4388214.5
$ws.Range("K71").Value = 21941072.5
$ws.Range("M71").Value = -21937328.5
$ws.Range("H113").Value = 2876.5557
$ws.Range("I113").Value = 3155.1072
$ws.Range("J113").Value = 1901.625
$ws.Range("K113").Value = 3155.1072
$ws.Range("L113").Value = 1901.625
$ws.Range("M113").Value = -985.1071999999999
$ws.Range("N113").Value = -6241.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2658.6
$ws.Range("J81").Value = 1996.75
$ws.Range("L81").Value = 3993.5
$ws.Range("N81").Value = -6115.5
$ws.Range("H84").Value = 2658.6
$ws.Range("J84").Value = 1996.75
$ws.Range("L84").Value = 19967.5
$ws.Range("N84").Value = -30575.5
$ws.Range("H107").Value = 339
$ws.Range("I107").Value = 326.8
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 980.4000000000001
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 939.5999999999999
$ws.Range("N107").Value = -5040
$ws.Range("H113").Value = 719.5454999999999
$ws.Range("I113").Value = 574.5625
$ws.Range("J113").Value = 1106.1666
$ws.Range("K113").Value = 1723.6875
$ws.Range("L113").Value = 3318.4998
$ws.Range("M113").Value = 446.3125
$ws.Range("N113").Value = -7658.4998
